$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 8443.866
$ws.Range("J17").Value = 8939.857
$ws.Range("L17").Value = 26819.571
$ws.Range("N17").Value = -27155.571
$ws.Range("H62").Value = 5240.8823
$ws.Range("I62").Value = 2735.3635
$ws.Range("K62").Value = 2735.3635
$ws.Range("M62").Value = -2111.3635
$ws.Range("H65").Value = 5240.8823
$ws.Range("I65").Value = 2735.3635
$ws.Range("K65").Value = 13676.8175
$ws.Range("M65").Value = -10556.8175
$ws.Range("H92").Value = 2807.64
$ws.Range("I92").Value = 1833.2858
$ws.Range("J92").Value = 4047.7273
$ws.Range("K92").Value = 1833.2858
$ws.Range("L92").Value = 4047.7273
$ws.Range("M92").Value = -585.2858000000001
$ws.Range("N92").Value = -6543.7273
$ws.Range("H115").Value = 541.625
$ws.Range("I115").Value = 333.2857
$ws.Range("K115").Value = 999.8571000000001
$ws.Range("M115").Value = 567.1428999999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2958.625
$ws.Range("I45").Value = 3445.8333
$ws.Range("K45").Value = 3445.8333
$ws.Range("M45").Value = -3068.8333
$ws.Range("H122").Value = 3194.1072
$ws.Range("I122").Value = 2649.0625
$ws.Range("K122").Value = 7947.1875
$ws.Range("M122").Value = -5497.1875
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 7235.5557
$ws.Range("J86").Value = 8466.666999999999
$ws.Range("L86").Value = 8466.666999999999
$ws.Range("N86").Value = -10712.667
$ws.Range("H89").Value = 7235.5557
$ws.Range("J89").Value = 8466.666999999999
$ws.Range("L89").Value = 42333.335
$ws.Range("N89").Value = -53565.335
$ws.Range("H107").Value = 1707.2941
$ws.Range("I107").Value = 1479.1538
$ws.Range("K107").Value = 1479.1538
$ws.Range("M107").Value = 440.8462
$ws.Range("H132").Value = 3738.1853
$ws.Range("I132").Value = 2674.4211
$ws.Range("K132").Value = 8023.263300000001
$ws.Range("M132").Value = -5493.263300000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 8268.333000000001
$ws.Range("I63").Value = 2895.5
$ws.Range("J63").Value = 19014
$ws.Range("K63").Value = 8686.5
$ws.Range("L63").Value = 57042
$ws.Range("M63").Value = -7937.5
$ws.Range("N63").Value = -58540
$ws.Range("H64").Value = 62504544
$ws.Range("I64").Value = 250001010
$ws.Range("J64").Value = 5720.1665
$ws.Range("K64").Value = 750003030
$ws.Range("L64").Value = 17160.4995
$ws.Range("M64").Value = -750002760
$ws.Range("N64").Value = -17700.4995
$ws.Range("H66").Value = 8268.333000000001
$ws.Range("I66").Value = 2895.5
$ws.Range("J66").Value = 19014
$ws.Range("K66").Value = 26059.5
$ws.Range("L66").Value = 171126
$ws.Range("M66").Value = -22315.5
$ws.Range("N66").Value = -178614
$ws.Range("H67").Value = 62504544
$ws.Range("I67").Value = 250001010
$ws.Range("J67").Value = 5720.1665
$ws.Range("K67").Value = 750003030
$ws.Range("L67").Value = 17160.4995
$ws.Range("M67").Value = -750002094
$ws.Range("N67").Value = -19032.4995
$ws.Range("H70").Value = 8274.666999999999
$ws.Range("I70").Value = 4912
$ws.Range("K70").Value = 14736
$ws.Range("M70").Value = -14421
$ws.Range("H73").Value = 8274.666999999999
$ws.Range("I73").Value = 4912
$ws.Range("K73").Value = 14736
$ws.Range("M73").Value = -13644
$ws.Range("H87").Value = 13249.5
$ws.Range("I87").Value = 13249.5
$ws.Range("K87").Value = 39748.5
$ws.Range("M87").Value = -38500.5
$ws.Range("H90").Value = 13249.5
$ws.Range("I90").Value = 13249.5
$ws.Range("K90").Value = 119245.5
$ws.Range("M90").Value = -113005.5
$ws.Range("H103").Value = 1216.2858
$ws.Range("I103").Value = 284.66666
$ws.Range("J103").Value = 1915
$ws.Range("K103").Value = 853.9999799999999
$ws.Range("L103").Value = 5745
$ws.Range("M103").Value = 25.00002000000006
$ws.Range("N103").Value = -7503
$ws.Range("H124").Value = 3059.4
$ws.Range("I124").Value = 2164.5
$ws.Range("J124").Value = 3283.125
$ws.Range("K124").Value = 6493.5
$ws.Range("L124").Value = 9849.375
$ws.Range("M124").Value = -1583.5
$ws.Range("N124").Value = -19669.375
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2311.4912
$ws.Range("I102").Value = 1631.8372
$ws.Range("K102").Value = 1631.8372
$ws.Range("M102").Value = -9.837199999999939
$ws.Range("H122").Value = 11695.421
$ws.Range("I122").Value = 16787.555
$ws.Range("J122").Value = 7112.5
$ws.Range("K122").Value = 50362.665
$ws.Range("L122").Value = 21337.5
$ws.Range("M122").Value = -47912.665
$ws.Range("N122").Value = -26237.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 165261.95
$ws.Range("I122").Value = 309748.62
$ws.Range("J122").Value = 8734.75
$ws.Range("K122").Value = 929245.86
$ws.Range("L122").Value = 26204.25
$ws.Range("M122").Value = -926795.86
$ws.Range("N122").Value = -31104.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 66692076
$ws.Range("I64").Value = 200000000
$ws.Range("K64").Value = 200000000
$ws.Range("M64").Value = -199999752
$ws.Range("H67").Value = 66692076
$ws.Range("I67").Value = 200000000
$ws.Range("K67").Value = 200000000
$ws.Range("M67").Value = -199999142
$ws.Range("H81").Value = 6366.778
$ws.Range("I81").Value = 3999
$ws.Range("J81").Value = 6662.75
$ws.Range("K81").Value = 7998
$ws.Range("L81").Value = 13325.5
$ws.Range("M81").Value = -6937
$ws.Range("N81").Value = -15447.5
$ws.Range("H84").Value = 6366.778
$ws.Range("I84").Value = 3999
$ws.Range("J84").Value = 6662.75
$ws.Range("K84").Value = 39990
$ws.Range("L84").Value = 66627.5
$ws.Range("M84").Value = -34686
$ws.Range("N84").Value = -77235.5
$ws.Range("H86").Value = 49000
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 49000
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H126").Value = 1340.25
$ws.Range("I126").Value = 1376.4615
$ws.Range("K126").Value = 4129.3845
$ws.Range("M126").Value = -1659.3845
$ws.Range("H136").Value = 4591.087
$ws.Range("I136").Value = 2053.1538
$ws.Range("K136").Value = 6159.4614
$ws.Range("M136").Value = -3609.4614
